$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A31").Value = 0.673673
$ws.Range("B31").Value = 0.892892
$ws.Range("C31").Value = 0.5116643051909798
$ws.Range("D31").Value = "query"
